# "Refined metadata to be additional tab"
# - Adds a new "metadata" worksheet after the existing "data" sheet, containing
#   a single header row + one data row describing the PanelApp query.
# - Refreshes the "time_taken" timestamps (column F) on the "data" sheet to the
#   time of the (re-)query.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update the per-gene query timestamps on the "data" sheet -----------
$newTimes = @(
    "2021-10-05 14:33:31.341442",
    "2021-10-05 14:33:31.341450",
    "2021-10-05 14:33:31.341453",
    "2021-10-05 14:33:31.341456",
    "2021-10-05 14:33:31.341458",
    "2021-10-05 14:33:31.341461",
    "2021-10-05 14:33:31.341464",
    "2021-10-05 14:33:31.341467",
    "2021-10-05 14:33:31.341470",
    "2021-10-05 14:33:31.341473",
    "2021-10-05 14:33:31.341476",
    "2021-10-05 14:33:31.341478",
    "2021-10-05 14:33:31.341481",
    "2021-10-05 14:33:31.341483",
    "2021-10-05 14:33:31.341486",
    "2021-10-05 14:33:31.341489",
    "2021-10-05 14:33:31.341491",
    "2021-10-05 14:33:31.341494",
    "2021-10-05 14:33:31.341497",
    "2021-10-05 14:33:31.341499",
    "2021-10-05 14:33:31.341502",
    "2021-10-05 14:33:31.341504",
    "2021-10-05 14:33:31.341507",
    "2021-10-05 14:33:31.341509",
    "2021-10-05 14:33:31.341512",
    "2021-10-05 14:33:31.341515",
    "2021-10-05 14:33:31.341518",
    "2021-10-05 14:33:31.341521",
    "2021-10-05 14:33:31.341524",
    "2021-10-05 14:33:31.341527",
    "2021-10-05 14:33:31.341529",
    "2021-10-05 14:33:31.341532",
    "2021-10-05 14:33:31.341535",
    "2021-10-05 14:33:31.341537",
    "2021-10-05 14:33:31.341540",
    "2021-10-05 14:33:31.341543",
    "2021-10-05 14:33:31.341545",
    "2021-10-05 14:33:31.341548"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $newTimes[$i]
}

# --- 2. Add the new "metadata" sheet, placed right after "data" ------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold / bordered / centered, matching the "data" sheet header
# style) - copy the format from the "data" header row so the same style
# record is reused rather than minting a near-duplicate one. The "data"
# header only spans B1:F1, so G1's format is seeded from F1 individually.
$dataSheet.Range("B1:F1").Copy() | Out-Null
$metaSheet.Range("B1:F1").PasteSpecial(-4122) | Out-Null
$dataSheet.Range("F1").Copy() | Out-Null
$metaSheet.Range("G1").PasteSpecial(-4122) | Out-Null

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (the pandas-style integer index in column A picks up the same
# style as the "data" sheet's index column).
$dataSheet.Range("A2").Copy() | Out-Null
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$metaSheet.Range("A2").Value = 0

$metaSheet.Range("B2").Value = "Common Variable Immunodeficiency"
$metaSheet.Range("C2").Value = 225

# "1.0" must round-trip as literal text (not get normalised to the number 1).
# Build it via a TEXT() formula in a scratch cell and paste-values it in, which
# yields a genuine text cell without leaving a stray number-format style behind.
$metaSheet.Range("Z1").Formula = '=TEXT(1,"0.0")'
$metaSheet.Range("Z1").Copy() | Out-Null
$metaSheet.Range("D2").PasteSpecial(-4163) | Out-Null
$metaSheet.Range("Z1").Clear() | Out-Null

$metaSheet.Range("E2").Value = "2021-08-09T02:01:56.409388Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:31.338189"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/225/?format=json"

$dataSheet.Activate()

Write-Output "metadata sheet added; data timestamps refreshed"
